$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated p-values (column F, "Pr(>F)") rounded to 3 significant/decimal digits
$ws.Range("F2").Value = 0.555
$ws.Range("F3").Value = 0.0114
$ws.Range("F4").Value = 0.577
$ws.Range("F7").Value = 0.1574
$ws.Range("F8").Value = 0.0854
$ws.Range("F9").Value = 0.9872
$ws.Range("F12").Value = 0.7038
$ws.Range("F13").Value = 0.7262
$ws.Range("F14").Value = 0.5914
$ws.Range("F17").Value = 0.0441
$ws.Range("F18").Value = 0.8305
$ws.Range("F19").Value = 0.351
$ws.Range("F22").Value = 0.3975
$ws.Range("F23").Value = 0.3941
$ws.Range("F24").Value = 0.1848
$ws.Range("F27").Value = 0.542
$ws.Range("F28").Value = 0.0098
$ws.Range("F29").Value = 0.5772
$ws.Range("F32").Value = 0.2036
$ws.Range("F33").Value = 0.6316
$ws.Range("F34").Value = 0.9873
